$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("B16").Value = 0.62136999999999998
$ws.Activate()
